$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 369.22223  # ALC H4: 305.66666 -> 369.22223
$ws.Cells.Item(4, 9).Value = 358.83334  # ALC I4: 267 -> 358.83334
$ws.Cells.Item(4, 10).Value = 390  # ALC J4: 363.66666 -> 390
$ws.Cells.Item(4, 11).Value = 358.83334  # ALC K4: 267 -> 358.83334
$ws.Cells.Item(4, 12).Value = 390  # ALC L4: 363.66666 -> 390
$ws.Cells.Item(4, 13).Value = -244.83334  # ALC M4: -153 -> -244.83334
$ws.Cells.Item(4, 14).Value = -618  # ALC N4: -591.66666 -> -618

$ws.Cells.Item(64, 8).Value = 3337.182  # ALC H64: 3579.875 -> 3337.182
$ws.Cells.Item(64, 9).Value = 2907  # ALC I64: 3341 -> 2907
$ws.Cells.Item(64, 10).Value = 3853.4  # ALC J64: 3659.5 -> 3853.4
$ws.Cells.Item(64, 11).Value = 2907  # ALC K64: 3341 -> 2907
$ws.Cells.Item(64, 12).Value = 3853.4  # ALC L64: 3659.5 -> 3853.4
$ws.Cells.Item(64, 13).Value = -2659  # ALC M64: -3093 -> -2659
$ws.Cells.Item(64, 14).Value = -4349.4  # ALC N64: -4155.5 -> -4349.4

$ws.Cells.Item(67, 8).Value = 3337.182  # ALC H67: 3579.875 -> 3337.182
$ws.Cells.Item(67, 9).Value = 2907  # ALC I67: 3341 -> 2907
$ws.Cells.Item(67, 10).Value = 3853.4  # ALC J67: 3659.5 -> 3853.4
$ws.Cells.Item(67, 11).Value = 2907  # ALC K67: 3341 -> 2907
$ws.Cells.Item(67, 12).Value = 3853.4  # ALC L67: 3659.5 -> 3853.4
$ws.Cells.Item(67, 13).Value = -2049  # ALC M67: -2483 -> -2049
$ws.Cells.Item(67, 14).Value = -5569.4  # ALC N67: -5375.5 -> -5569.4

$ws.Cells.Item(76, 8).Value = 2182238.2  # ALC H76: 2182256 -> 2182238.2
$ws.Cells.Item(76, 9).Value = 2318281  # ALC I76: 2318300 -> 2318281
$ws.Cells.Item(76, 11).Value = 2318281  # ALC K76: 2318300 -> 2318281
$ws.Cells.Item(76, 13).Value = -2317966  # ALC M76: -2317985 -> -2317966

$ws.Cells.Item(79, 8).Value = 2182238.2  # ALC H79: 2182256 -> 2182238.2
$ws.Cells.Item(79, 9).Value = 2318281  # ALC I79: 2318300 -> 2318281
$ws.Cells.Item(79, 11).Value = 2318281  # ALC K79: 2318300 -> 2318281
$ws.Cells.Item(79, 13).Value = -2317189  # ALC M79: -2317208 -> -2317189

$ws.Cells.Item(80, 8).Value = 2860.08  # ALC H80: 2971.3333 -> 2860.08
$ws.Cells.Item(80, 9).Value = 884.7646999999999  # ALC I80: 925.0625 -> 884.7646999999999
$ws.Cells.Item(80, 10).Value = 7057.625  # ALC J80: 7063.875 -> 7057.625
$ws.Cells.Item(80, 11).Value = 2654.2941  # ALC K80: 2775.1875 -> 2654.2941
$ws.Cells.Item(80, 12).Value = 21172.875  # ALC L80: 21191.625 -> 21172.875
$ws.Cells.Item(80, 13).Value = -1656.2941  # ALC M80: -1777.1875 -> -1656.2941
$ws.Cells.Item(80, 14).Value = -23168.875  # ALC N80: -23187.625 -> -23168.875

$ws.Cells.Item(83, 8).Value = 2860.08  # ALC H83: 2971.3333 -> 2860.08
$ws.Cells.Item(83, 9).Value = 884.7646999999999  # ALC I83: 925.0625 -> 884.7646999999999
$ws.Cells.Item(83, 10).Value = 7057.625  # ALC J83: 7063.875 -> 7057.625
$ws.Cells.Item(83, 11).Value = 7962.882299999999  # ALC K83: 8325.5625 -> 7962.882299999999
$ws.Cells.Item(83, 12).Value = 63518.625  # ALC L83: 63574.875 -> 63518.625
$ws.Cells.Item(83, 13).Value = -2970.882299999999  # ALC M83: -3333.5625 -> -2970.882299999999
$ws.Cells.Item(83, 14).Value = -73502.625  # ALC N83: -73558.875 -> -73502.625

$ws.Cells.Item(87, 8).Value = 10975.553  # ALC H87: 10982.361 -> 10975.553
$ws.Cells.Item(87, 9).Value = 2910.5  # ALC I87: 3410.5 -> 2910.5
$ws.Cells.Item(87, 10).Value = 11423.611  # ALC J87: 11427.765 -> 11423.611
$ws.Cells.Item(87, 11).Value = 2910.5  # ALC K87: 3410.5 -> 2910.5
$ws.Cells.Item(87, 12).Value = 11423.611  # ALC L87: 11427.765 -> 11423.611
$ws.Cells.Item(87, 13).Value = -1662.5  # ALC M87: -2162.5 -> -1662.5
$ws.Cells.Item(87, 14).Value = -13919.611  # ALC N87: -13923.765 -> -13919.611

$ws.Cells.Item(90, 8).Value = 10975.553  # ALC H90: 10982.361 -> 10975.553
$ws.Cells.Item(90, 9).Value = 2910.5  # ALC I90: 3410.5 -> 2910.5
$ws.Cells.Item(90, 10).Value = 11423.611  # ALC J90: 11427.765 -> 11423.611
$ws.Cells.Item(90, 11).Value = 8731.5  # ALC K90: 10231.5 -> 8731.5
$ws.Cells.Item(90, 12).Value = 34270.833  # ALC L90: 34283.295 -> 34270.833
$ws.Cells.Item(90, 13).Value = -2491.5  # ALC M90: -3991.5 -> -2491.5
$ws.Cells.Item(90, 14).Value = -46750.833  # ALC N90: -46763.295 -> -46750.833

$ws.Cells.Item(136, 8).Value = 0  # ALC H136: 48000 -> 0
$ws.Cells.Item(136, 10).Value = 0  # ALC J136: 48000 -> 0
$ws.Cells.Item(136, 12).Value = 0  # ALC L136: 48000 -> 0
$ws.Cells.Item(136, 14).ClearContents()  # ALC N136: -58200 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3832.9443  # ARM H61: 2751.8386 -> 3832.9443
$ws.Cells.Item(61, 9).Value = 4291.923  # ARM I61: 2729.5417 -> 4291.923
$ws.Cells.Item(61, 10).Value = 2639.6  # ARM J61: 2828.2856 -> 2639.6
$ws.Cells.Item(61, 11).Value = 4291.923  # ARM K61: 2729.5417 -> 4291.923
$ws.Cells.Item(61, 12).Value = 2639.6  # ARM L61: 2828.2856 -> 2639.6
$ws.Cells.Item(61, 13).Value = -4079.923  # ARM M61: -2517.5417 -> -4079.923
$ws.Cells.Item(61, 14).Value = -3063.6  # ARM N61: -3252.2856 -> -3063.6

$ws.Cells.Item(122, 8).Value = 21047.139  # ARM H122: 5084.724 -> 21047.139
$ws.Cells.Item(122, 9).Value = 25904.143  # ARM I122: 5303.4546 -> 25904.143
$ws.Cells.Item(122, 10).Value = 4047.625  # ARM J122: 4397.2856 -> 4047.625
$ws.Cells.Item(122, 11).Value = 77712.429  # ARM K122: 15910.3638 -> 77712.429
$ws.Cells.Item(122, 12).Value = 12142.875  # ARM L122: 13191.8568 -> 12142.875
$ws.Cells.Item(122, 13).Value = -75262.429  # ARM M122: -13460.3638 -> -75262.429
$ws.Cells.Item(122, 14).Value = -17042.875  # ARM N122: -18091.8568 -> -17042.875

$ws.Cells.Item(130, 8).Value = 56594  # ARM H130: 60525 -> 56594
$ws.Cells.Item(130, 10).Value = 56594  # ARM J130: 60525 -> 56594
$ws.Cells.Item(130, 12).Value = 56594  # ARM L130: 60525 -> 56594
$ws.Cells.Item(130, 14).Value = -66634  # ARM N130: -70565 -> -66634

$ws.Cells.Item(132, 8).Value = 2225.3096  # ARM H132: 2126.6 -> 2225.3096
$ws.Cells.Item(132, 9).Value = 2054.1  # ARM I132: 1964.2812 -> 2054.1
$ws.Cells.Item(132, 10).Value = 2653.3333  # ARM J132: 2526.1538 -> 2653.3333
$ws.Cells.Item(132, 11).Value = 6162.299999999999  # ARM K132: 5892.8436 -> 6162.299999999999
$ws.Cells.Item(132, 12).Value = 7959.999899999999  # ARM L132: 7578.4614 -> 7959.999899999999
$ws.Cells.Item(132, 13).Value = -3632.299999999999  # ARM M132: -3362.8436 -> -3632.299999999999
$ws.Cells.Item(132, 14).Value = -13019.9999  # ARM N132: -12638.4614 -> -13019.9999

$ws.Cells.Item(136, 8).Value = 3832.9443  # ARM H136: 2751.8386 -> 3832.9443
$ws.Cells.Item(136, 9).Value = 4291.923  # ARM I136: 2729.5417 -> 4291.923
$ws.Cells.Item(136, 10).Value = 2639.6  # ARM J136: 2828.2856 -> 2639.6
$ws.Cells.Item(136, 11).Value = 12875.769  # ARM K136: 8188.625100000001 -> 12875.769
$ws.Cells.Item(136, 12).Value = 7918.799999999999  # ARM L136: 8484.856800000001 -> 7918.799999999999
$ws.Cells.Item(136, 13).Value = -10325.769  # ARM M136: -5638.625100000001 -> -10325.769
$ws.Cells.Item(136, 14).Value = -13018.8  # ARM N136: -13584.8568 -> -13018.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4111.75  # BSM H107: 4168.647 -> 4111.75
$ws.Cells.Item(107, 9).Value = 4314.5938  # BSM I107: 4243.0625 -> 4314.5938
$ws.Cells.Item(107, 10).Value = 2489  # BSM J107: 2978 -> 2489
$ws.Cells.Item(107, 11).Value = 4314.5938  # BSM K107: 4243.0625 -> 4314.5938
$ws.Cells.Item(107, 12).Value = 2489  # BSM L107: 2978 -> 2489
$ws.Cells.Item(107, 13).Value = -2394.5938  # BSM M107: -2323.0625 -> -2394.5938
$ws.Cells.Item(107, 14).Value = -6329  # BSM N107: -6818 -> -6329

$ws.Cells.Item(134, 8).Value = 10173.667  # BSM H134: 5453.1113 -> 10173.667
$ws.Cells.Item(134, 9).Value = 1051  # BSM I134: 853.26666 -> 1051
$ws.Cells.Item(134, 10).Value = 28419  # BSM J134: 28452.334 -> 28419
$ws.Cells.Item(134, 11).Value = 3153  # BSM K134: 2559.79998 -> 3153
$ws.Cells.Item(134, 12).Value = 85257  # BSM L134: 85357.00199999999 -> 85257
$ws.Cells.Item(134, 13).Value = -618  # BSM M134: -24.79997999999978 -> -618
$ws.Cells.Item(134, 14).Value = -90327  # BSM N134: -90427.00199999999 -> -90327

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1303.9412  # CRP H58: 1221.8948 -> 1303.9412
$ws.Cells.Item(58, 9).Value = 797.375  # CRP I58: 738 -> 797.375
$ws.Cells.Item(58, 10).Value = 1754.2222  # CRP J58: 2051.4285 -> 1754.2222
$ws.Cells.Item(58, 11).Value = 797.375  # CRP K58: 738 -> 797.375
$ws.Cells.Item(58, 12).Value = 1754.2222  # CRP L58: 2051.4285 -> 1754.2222
$ws.Cells.Item(58, 13).Value = -594.375  # CRP M58: -535 -> -594.375
$ws.Cells.Item(58, 14).Value = -2160.2222  # CRP N58: -2457.4285 -> -2160.2222

$ws.Cells.Item(136, 8).Value = 1303.9412  # CRP H136: 1221.8948 -> 1303.9412
$ws.Cells.Item(136, 9).Value = 797.375  # CRP I136: 738 -> 797.375
$ws.Cells.Item(136, 10).Value = 1754.2222  # CRP J136: 2051.4285 -> 1754.2222
$ws.Cells.Item(136, 11).Value = 2392.125  # CRP K136: 2214 -> 2392.125
$ws.Cells.Item(136, 12).Value = 5262.6666  # CRP L136: 6154.2855 -> 5262.6666
$ws.Cells.Item(136, 13).Value = 157.875  # CRP M136: 336 -> 157.875
$ws.Cells.Item(136, 14).Value = -10362.6666  # CRP N136: -11254.2855 -> -10362.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 989.10205  # CUL H113: 999.2917 -> 989.10205
$ws.Cells.Item(113, 9).Value = 733.3333  # CUL I113: 850 -> 733.3333
$ws.Cells.Item(113, 11).Value = 2199.9999  # CUL K113: 2550 -> 2199.9999
$ws.Cells.Item(113, 13).Value = -29.9998999999998  # CUL M113: -380 -> -29.9998999999998

$ws.Cells.Item(131, 8).Value = 5747990  # CUL H131: 5320017.5 -> 5747990
$ws.Cells.Item(131, 9).Value = 868.4  # CUL I131: 869.35 -> 868.4
$ws.Cells.Item(131, 10).Value = 7463548.5  # CUL J131: 6757625 -> 7463548.5
$ws.Cells.Item(131, 11).Value = 2605.2  # CUL K131: 2608.05 -> 2605.2
$ws.Cells.Item(131, 12).Value = 22390645.5  # CUL L131: 20272875 -> 22390645.5
$ws.Cells.Item(131, 13).Value = 2434.8  # CUL M131: 2431.95 -> 2434.8
$ws.Cells.Item(131, 14).Value = -22400725.5  # CUL N131: -20282955 -> -22400725.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 43511.184  # GSM H80: 64517.055 -> 43511.184
$ws.Cells.Item(80, 9).Value = 71581.375  # GSM I80: 103527.91 -> 71581.375
$ws.Cells.Item(80, 10).Value = 2681.818  # GSM J80: 3214.2856 -> 2681.818
$ws.Cells.Item(80, 11).Value = 71581.375  # GSM K80: 103527.91 -> 71581.375
$ws.Cells.Item(80, 12).Value = 2681.818  # GSM L80: 3214.2856 -> 2681.818
$ws.Cells.Item(80, 13).Value = -70583.375  # GSM M80: -102529.91 -> -70583.375
$ws.Cells.Item(80, 14).Value = -4677.818  # GSM N80: -5210.2856 -> -4677.818

$ws.Cells.Item(83, 8).Value = 43511.184  # GSM H83: 64517.055 -> 43511.184
$ws.Cells.Item(83, 9).Value = 71581.375  # GSM I83: 103527.91 -> 71581.375
$ws.Cells.Item(83, 10).Value = 2681.818  # GSM J83: 3214.2856 -> 2681.818
$ws.Cells.Item(83, 11).Value = 357906.875  # GSM K83: 517639.55 -> 357906.875
$ws.Cells.Item(83, 12).Value = 13409.09  # GSM L83: 16071.428 -> 13409.09
$ws.Cells.Item(83, 13).Value = -352914.875  # GSM M83: -512647.55 -> -352914.875
$ws.Cells.Item(83, 14).Value = -23393.09  # GSM N83: -26055.428 -> -23393.09

$ws.Cells.Item(111, 8).Value = 0  # GSM H111: 19380 -> 0
$ws.Cells.Item(111, 10).Value = 0  # GSM J111: 19380 -> 0
$ws.Cells.Item(111, 12).Value = 0  # GSM L111: 19380 -> 0
$ws.Cells.Item(111, 14).ClearContents()  # GSM N111: -25514 -> (removed)

$ws.Cells.Item(122, 8).Value = 2143.9534  # GSM H122: 2692.077 -> 2143.9534
$ws.Cells.Item(122, 9).Value = 2013.95  # GSM I122: 2318.1428 -> 2013.95
$ws.Cells.Item(122, 10).Value = 2257  # GSM J122: 3128.3333 -> 2257
$ws.Cells.Item(122, 11).Value = 6041.85  # GSM K122: 6954.428400000001 -> 6041.85
$ws.Cells.Item(122, 12).Value = 6771  # GSM L122: 9384.999899999999 -> 6771
$ws.Cells.Item(122, 13).Value = -3591.85  # GSM M122: -4504.428400000001 -> -3591.85
$ws.Cells.Item(122, 14).Value = -11671  # GSM N122: -14284.9999 -> -11671

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3208  # LTW H93: 3622.3076 -> 3208
$ws.Cells.Item(93, 9).Value = 3015.5715  # LTW I93: 4066.5 -> 3015.5715
$ws.Cells.Item(93, 10).Value = 3357.6667  # LTW J93: 3424.889 -> 3357.6667
$ws.Cells.Item(93, 11).Value = 3015.5715  # LTW K93: 4066.5 -> 3015.5715
$ws.Cells.Item(93, 12).Value = 3357.6667  # LTW L93: 3424.889 -> 3357.6667
$ws.Cells.Item(93, 13).Value = -1767.5715  # LTW M93: -2818.5 -> -1767.5715
$ws.Cells.Item(93, 14).Value = -5853.6667  # LTW N93: -5920.889 -> -5853.6667

$ws.Cells.Item(110, 8).Value = 35762.668  # LTW H110: 39021 -> 35762.668
$ws.Cells.Item(110, 10).Value = 35762.668  # LTW J110: 39021 -> 35762.668
$ws.Cells.Item(110, 12).Value = 35762.668  # LTW L110: 39021 -> 35762.668
$ws.Cells.Item(110, 14).Value = -43942.668  # LTW N110: -47201 -> -43942.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6208.4  # WVR H62: 5692.5386 -> 6208.4
$ws.Cells.Item(62, 9).Value = 5734.7144  # WVR I62: 6500 -> 5734.7144
$ws.Cells.Item(62, 10).Value = 6622.875  # WVR J62: 5000.4287 -> 6622.875
$ws.Cells.Item(62, 11).Value = 5734.7144  # WVR K62: 6500 -> 5734.7144
$ws.Cells.Item(62, 12).Value = 6622.875  # WVR L62: 5000.4287 -> 6622.875
$ws.Cells.Item(62, 13).Value = -5110.7144  # WVR M62: -5876 -> -5110.7144
$ws.Cells.Item(62, 14).Value = -7870.875  # WVR N62: -6248.4287 -> -7870.875

$ws.Cells.Item(65, 8).Value = 6208.4  # WVR H65: 5692.5386 -> 6208.4
$ws.Cells.Item(65, 9).Value = 5734.7144  # WVR I65: 6500 -> 5734.7144
$ws.Cells.Item(65, 10).Value = 6622.875  # WVR J65: 5000.4287 -> 6622.875
$ws.Cells.Item(65, 11).Value = 28673.572  # WVR K65: 32500 -> 28673.572
$ws.Cells.Item(65, 12).Value = 33114.375  # WVR L65: 25002.1435 -> 33114.375
$ws.Cells.Item(65, 13).Value = -25553.572  # WVR M65: -29380 -> -25553.572
$ws.Cells.Item(65, 14).Value = -39354.375  # WVR N65: -31242.1435 -> -39354.375
